# Término da atividade 'Prova 2des' - LIMA
#
# Update the age-bracket labels in the "FasesDaVida" lookup table (K2:N2)
# to add spacing around the comparison operators, and move the active
# selection from I5 to N3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "> 0 <= 13"
$ws.Range("L2").Value = ">13 <= 20"
$ws.Range("M2").Value = "> 20 <= 60"
$ws.Range("N2").Value = "> 60"

$ws.Range("N3").Select()
